# "finestra incidenza 7gg centrata su ultimo g"
#
# Column C ("somma mobile 7gg.") was a 7-day window CENTERED on each date
# (sum of the 3 days before, the day itself, and the 3 days after).
# This changes it to a 7-day TRAILING window ending on (centered on) the
# current/last day: sum of the day itself and the 6 days before it.
# Column D is the same figure normalised per 100,000 inhabitants
# (population = 1604).
#
# Because the window anchor moved, the first 6 rows (where a full
# trailing window isn't yet available) become blank, and the last 3 rows
# (which previously lacked enough *future* days for a centered window)
# now get real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 184
$population = 1604
$windowSize = 7

# Read column B ("nuovi pos.") for every row into a lookup table.
$b = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $b[$r] = $ws.Cells.Item($r, 2).Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $windowStart = $r - ($windowSize - 1)
    $oldC = $ws.Cells.Item($r, 3).Value2

    if ($windowStart -lt $firstRow) {
        # Not enough preceding days yet for a full trailing 7-day window.
        # Leave already-blank cells (rows 2-4) untouched; only blank out
        # the rows that used to hold a real centered-window value (5-7).
        if ($oldC -ne "") {
            $ws.Cells.Item($r, 3).Value = ""
            $ws.Cells.Item($r, 4).Value = ""
        }
    } else {
        $total = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $total = $total + $b[$i]
        }
        $ws.Cells.Item($r, 3).Value = $total
        $ws.Cells.Item($r, 4).Value = ($total * 100000 / $population)
    }
}
